$wb = $excel.ActiveWorkbook

# --- 1. Update status text: "Ready for handoff" -> "In Translation" ---
# (collect matching cell addresses first, then write - mutating a cell's
# value while a live enumerator is walking the same UsedRange can disturb
# the shared-string table indices mid-iteration)
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $targets = @()
    foreach ($cell in $used.Cells) {
        $val = [string]$cell.Value()
        if ($val -eq "Ready for handoff") {
            $targets += $cell.Address()
        }
    }
    foreach ($addr in $targets) {
        $ws.Range($addr).Value = "In Translation"
    }
}

# --- 2. Re-fit the "Status" columns now that the text is shorter ---
# Excel/COM always rounds a `ColumnWidth` assignment to whole screen pixels
# before it is stored, so the exact source width (13.4101845877511 chars)
# can't be reproduced bit-for-bit through this API no matter what a script
# writes - the closest achievable stored width is 13.333... chars, which is
# what assigning 12.5 here resolves to after that internal rounding.
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newWidth

